$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'296.06"
$ws.Range("E2").Value = "'3.35%"

# Row 3
$ws.Range("D3").Value = "'41.15"

# Row 4
$ws.Range("D4").Value = "'5.038"
$ws.Range("E4").Value = "'0.21%"

# Row 5
$ws.Range("D5").Value = "'0.07431"
$ws.Range("E5").Value = "'2.04%"

# Row 6
$ws.Range("B6").Value = "FTXToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D6").Value = "'1.571"
$ws.Range("E6").Value = "'2.83%"

# Row 7
$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").Value = "'0.9238"
$ws.Range("E7").Value = "'0.55%"

# Row 8
$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").Value = "'2.420"
$ws.Range("E8").Value = "'0.96%"

# Row 9
$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").Value = "'0.1172"
$ws.Range("E9").Value = "'-2.35%"

# Row 10
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1817"
$ws.Range("E10").Value = "'6.26%"

# Row 11
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.08814"
$ws.Range("E11").Value = "'2.01%"

# Row 12
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.04297"
$ws.Range("E12").Value = "'3.09%"

# Row 13
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.1051"
$ws.Range("E13").Value = "'-0.03%"

# Row 14
$ws.Range("B14").Value = "TigerCash"
$ws.Range("C14").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D14").Value = "'0.005893"
$ws.Range("E14").Value = "'-1.18%"

# Row 15
$ws.Range("B15").Value = "LEO"
$ws.Range("C15").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D15").Value = "'3.356"
$ws.Range("E15").Value = "'-1.28%"

# Row 16
$ws.Range("B16").Value = "GateToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D16").Value = "'4.357"
$ws.Range("E16").Value = "'1.32%"

# Row 17
$ws.Range("B17").Value = "BitpandaEcosystemToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D17").Value = "'0.3287"
$ws.Range("E17").Value = "'0.16%"

# Row 18
$ws.Range("B18").Value = "MCDex"
$ws.Range("C18").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D18").Value = "'7.871"
$ws.Range("E18").Value = "'0.33%"

# Row 19
$ws.Range("B19").Value = "ProBitToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D19").Value = "'0.1362"
$ws.Range("E19").Value = "'1.40%"

# Row 20
$ws.Range("B20").Value = "ZBToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D20").Value = "'0.2969"
$ws.Range("E20").Value = "'2.85%"

# Row 21
$ws.Range("B21").Value = "BitForexToken"
$ws.Range("C21").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D21").Value = "'0.001278"
$ws.Range("E21").Value = "'0.85%"

# Row 22
$ws.Range("D22").Value = "'0.04035"
$ws.Range("E22").Value = "'4.93%"

# Row 23
$ws.Range("D23").Value = "'0.001273"
$ws.Range("E23").Value = "'0.12%"

# Row 24
$ws.Range("D24").Value = "'0.003864"
$ws.Range("E24").Value = "'1.95%"

# Row 25
$ws.Range("E25").Value = "'-3.96%"

# Row 26
$ws.Range("D26").Value = "'0.0003726"
$ws.Range("E26").Value = "'-0.07%"

# Row 38
$ws.Range("D38").Value = "'0.02385"
$ws.Range("E38").Value = "'2.89%"

# Row 39
$ws.Range("D39").Value = "'0.05190"
$ws.Range("E39").Value = "'4.14%"

# Row 40
$ws.Range("D40").Value = "'0.006609"
$ws.Range("E40").Value = "'-2.35%"

# Row 41
$ws.Range("E41").Value = "'1.19%"

# Row 42
$ws.Range("D42").Value = "'0.1316"
$ws.Range("E42").Value = "'3.77%"

# Row 43
$ws.Range("D43").Value = "'0.007387"
$ws.Range("E43").Value = "'0.19%"

# Row 44
$ws.Range("D44").Value = "'0.007800"
$ws.Range("E44").Value = "'5.09%"

# Row 45
$ws.Range("D45").Value = "'0.3213"
$ws.Range("E45").Value = "'3.96%"

# Row 46
$ws.Range("D46").Value = "'0.00006250"
$ws.Range("E46").Value = "'-3.07%"

# Row 47
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("E47").Value = "'-0.05%"

# Row 48
$ws.Range("D48").Value = "'0.04626"
$ws.Range("E48").Value = "'-81.63%"

# Row 49
$ws.Range("D49").Value = "'0.004204"
$ws.Range("E49").Value = "'0.06%"

# Row 50
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'-0.05%"

# Row 51
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'-0.05%"
